$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated regression results (fix for missing values bug)
$values = @{
    2  = @(0.3666, 0.579, 0.6331, 0.5285, -0.786, 1.5191)
    3  = @(-0.0254, 0.1011, -0.251, 0.8025, -0.2266, 0.1758)
    4  = @(0.0264, 0.1232, 0.2142, 0.8309, -0.2188, 0.2716)
    5  = @(0.3646, 0.0909, 4.01, 0.0001, 0.1836, 0.5455)
    6  = @(0.4022, 0.1059, 3.7972, 0.0003, 0.1914, 0.6131)
    7  = @(-0.031, 0.1018, -0.3046, 0.7615, -0.2336, 0.1716)
    8  = @(0.1112, 0.0783, 1.4201, 0.1595, -0.0446, 0.267)
    9  = @(-0.0932, 0.0811, -1.15, 0.2536, -0.2546, 0.0681)
    10 = @(0.0033, 0.0768, 0.0431, 0.9657, -0.1495, 0.1562)
    11 = @(-0.0345, 0.1926, -0.1793, 0.8582, -0.4179, 0.3488)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $col = 2 + $i  # Column B starts at index 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$i]
    }
}
